$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.036.33'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.564.17'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  +0.79%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('E6').Value = '  +0.42%  '
$ws.Range('E7').Value = '  +0.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.93'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.248'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0596'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0862'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = '1.785.68'
$ws.Range('E12').Value = '  +0.76%  '
$ws.Range('D13').Value = '1.571.79'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.518'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '27.004.31'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.93'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '0.0₃0703'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '215.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.105'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0472'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.03%  '
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.17'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.07%  '
$ws.Range('D34').Value = '1.430.39'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +16.27%  '
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.808'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.59%  '
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.28%  '
$ws.Range('D47').Value = '1.703.28'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.56'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('E50').Value = '  -0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0960'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.45%  '
